# Update lab annotation file names and lab_id values for several labs.
# Column A  = lab_file (the source filename for that lab's annotations)
# Column AB = lab_id   (the short lab identifier)
#
# Changes:
#   afekta: peak_evidence_rt_grouped_manual_fixed_afekta.xlsx -> fixed_annotation_afekta.xlsx
#           afekta -> _afekta
#   cembio: peak_evidence_rt_grouped_manual_fixed_cembio.xlsx -> fixed_annotation_cembio.xlsx
#           cembio -> _cembio
#   hmgu:   peak_evidence_rt_grouped_manual_fixed_hmgu.xlsx   -> fixed_annotation_hmgu.xlsx
#           hmgu -> _hmgu
#   icl:    peak_evidence_rt_grouped_manual_fixed_icl.xlsx    -> fixed_annotation_icl.xlsx
#           icl -> _icl

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colA = $ws.Columns.Item(1)
$colAB = $ws.Columns.Item(28)

# xlWhole = 1 (match the entire cell contents, not a substring)
$xlWhole = 1

$colA.Replace("peak_evidence_rt_grouped_manual_fixed_afekta.xlsx", "fixed_annotation_afekta.xlsx", $xlWhole)
$colAB.Replace("afekta", "_afekta", $xlWhole)

$colA.Replace("peak_evidence_rt_grouped_manual_fixed_cembio.xlsx", "fixed_annotation_cembio.xlsx", $xlWhole)
$colAB.Replace("cembio", "_cembio", $xlWhole)

$colA.Replace("peak_evidence_rt_grouped_manual_fixed_hmgu.xlsx", "fixed_annotation_hmgu.xlsx", $xlWhole)
$colAB.Replace("hmgu", "_hmgu", $xlWhole)

$colA.Replace("peak_evidence_rt_grouped_manual_fixed_icl.xlsx", "fixed_annotation_icl.xlsx", $xlWhole)
$colAB.Replace("icl", "_icl", $xlWhole)
